$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "bleu" status label to "noir"
$ws.Cells.Replace("bleu", "noir", 1)

# Correct the status_name wording
$ws.Cells.Replace("pas de résultat ni de publication", "pas de résultat postés ni publiés", 1)
$ws.Cells.Replace("résultat et / ou publication posté dans les 12 mois", "résultat postés ou publiés dans les 12 mois", 1)
$ws.Cells.Replace("résultat et / ou publication posté dans les 36 mois", "résultat postés ou publiés dans les 36 mois", 1)
